$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 7666.6665
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 7666.6665
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 7666.6665
$ws.Range("N21").Value = -8602.666499999999
$ws.Range("M21").ClearContents()

$ws.Range("H23").Value = 7666.6665
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 7666.6665
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 7666.6665
$ws.Range("N23").Value = -8134.6665
$ws.Range("M23").ClearContents()

$ws.Range("H34").Value = 7565.5713
$ws.Range("I34").Value = 2782
$ws.Range("K34").Value = 2782
$ws.Range("M34").Value = -2579

$ws.Range("H36").Value = 7565.5713
$ws.Range("I36").Value = 2782
$ws.Range("K36").Value = 2782
$ws.Range("M36").Value = -2067

$ws.Range("H74").Value = 5209.4814
$ws.Range("I74").Value = 7118.625
$ws.Range("J74").Value = 4405.6313
$ws.Range("K74").Value = 7118.625
$ws.Range("L74").Value = 4405.6313
$ws.Range("M74").Value = -6182.625
$ws.Range("N74").Value = -6277.6313

$ws.Range("H77").Value = 5209.4814
$ws.Range("I77").Value = 7118.625
$ws.Range("J77").Value = 4405.6313
$ws.Range("K77").Value = 35593.125
$ws.Range("L77").Value = 22028.1565
$ws.Range("M77").Value = -30913.125
$ws.Range("N77").Value = -31388.1565

$ws.Range("H92").Value = 2612.2273
$ws.Range("I92").Value = 2993.1052
$ws.Range("J92").Value = 200
$ws.Range("K92").Value = 2993.1052
$ws.Range("L92").Value = 200
$ws.Range("M92").Value = -1745.1052
$ws.Range("N92").Value = -2696

$ws.Range("H135").Value = 1662.3846
$ws.Range("I135").Value = 1212.375
$ws.Range("J135").Value = 2382.4
$ws.Range("K135").Value = 10911.375
$ws.Range("L135").Value = 21441.6
$ws.Range("M135").Value = -8376.375
$ws.Range("N135").Value = -26511.6

$ws.Range("H141").Value = 722.8611
$ws.Range("I141").Value = 647.7353000000001
$ws.Range("J141").Value = 2000
$ws.Range("K141").Value = 1943.2059
$ws.Range("L141").Value = 6000
$ws.Range("M141").Value = 3236.7941
$ws.Range("N141").Value = -16360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3749.06
$ws.Range("I32").Value = 2998.1929
$ws.Range("J32").Value = 7415.0586
$ws.Range("K32").Value = 2998.1929
$ws.Range("L32").Value = 7415.0586
$ws.Range("M32").Value = -2711.1929
$ws.Range("N32").Value = -7989.0586

$ws.Range("H122").Value = 3397.1592
$ws.Range("I122").Value = 3474.0967
$ws.Range("J122").Value = 3213.6924
$ws.Range("K122").Value = 10422.2901
$ws.Range("L122").Value = 9641.0772
$ws.Range("M122").Value = -7972.2901
$ws.Range("N122").Value = -14541.0772

$ws.Range("H132").Value = 2865.7324
$ws.Range("I132").Value = 2484.9075
$ws.Range("J132").Value = 4075.4119
$ws.Range("K132").Value = 7454.7225
$ws.Range("L132").Value = 12226.2357
$ws.Range("M132").Value = -4924.7225
$ws.Range("N132").Value = -17286.2357

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1186.2667
$ws.Range("I94").Value = 1054.5
$ws.Range("J94").Value = 1234.1818
$ws.Range("K94").Value = 1054.5
$ws.Range("L94").Value = 1234.1818
$ws.Range("M94").Value = -603.5
$ws.Range("N94").Value = -2136.1818

$ws.Range("H107").Value = 2985.5
$ws.Range("I107").Value = 2567
$ws.Range("J107").Value = 3311
$ws.Range("K107").Value = 2567
$ws.Range("L107").Value = 3311
$ws.Range("M107").Value = -647
$ws.Range("N107").Value = -7151

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2802.8774
$ws.Range("I31").Value = 1927.8108
$ws.Range("J31").Value = 5501
$ws.Range("K31").Value = 1927.8108
$ws.Range("L31").Value = 5501
$ws.Range("M31").Value = -1632.8108
$ws.Range("N31").Value = -6091

$ws.Range("H34").Value = 2802.8774
$ws.Range("I34").Value = 1927.8108
$ws.Range("J34").Value = 5501
$ws.Range("K34").Value = 1927.8108
$ws.Range("L34").Value = 5501
$ws.Range("M34").Value = -1725.8108
$ws.Range("N34").Value = -5905

$ws.Range("H132").Value = 2656.4375
$ws.Range("I132").Value = 1600.3636
$ws.Range("K132").Value = 4801.0908
$ws.Range("M132").Value = -2271.0908

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1054.1034
$ws.Range("I131").Value = 510
$ws.Range("J131").Value = 1083.7819
$ws.Range("K131").Value = 1530
$ws.Range("L131").Value = 3251.3457
$ws.Range("M131").Value = 3510
$ws.Range("N131").Value = -13331.3457

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2148.95
$ws.Range("I82").Value = 1508.25
$ws.Range("J82").Value = 3110
$ws.Range("K82").Value = 1508.25
$ws.Range("L82").Value = 3110
$ws.Range("M82").Value = -1147.25
$ws.Range("N82").Value = -3832

$ws.Range("H85").Value = 2148.95
$ws.Range("I85").Value = 1508.25
$ws.Range("J85").Value = 3110
$ws.Range("K85").Value = 1508.25
$ws.Range("L85").Value = 3110
$ws.Range("M85").Value = -260.25
$ws.Range("N85").Value = -5606

$ws.Range("H93").Value = 1406.1794
$ws.Range("I93").Value = 1423.7826
$ws.Range("J93").Value = 1380.875
$ws.Range("K93").Value = 1423.7826
$ws.Range("L93").Value = 1380.875
$ws.Range("M93").Value = -175.7826
$ws.Range("N93").Value = -3876.875

$ws.Range("H100").Value = 125007750
$ws.Range("I100").Value = 25450
$ws.Range("J100").Value = 166668510
$ws.Range("K100").Value = 25450
$ws.Range("L100").Value = 166668510
$ws.Range("M100").Value = -24909
$ws.Range("N100").Value = -166669592

$ws.Range("H132").Value = 13259.286
$ws.Range("I132").Value = 6050.5
$ws.Range("J132").Value = 16142.8
$ws.Range("K132").Value = 18151.5
$ws.Range("L132").Value = 48428.39999999999
$ws.Range("M132").Value = -15621.5
$ws.Range("N132").Value = -53488.39999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3138.2222
$ws.Range("I96").Value = 2760
$ws.Range("J96").Value = 3611
$ws.Range("K96").Value = 2760
$ws.Range("L96").Value = 3611
$ws.Range("M96").Value = -1387
$ws.Range("N96").Value = -6357

$ws.Range("H104").Value = 23624.428
$ws.Range("J104").Value = 23624.428
$ws.Range("L104").Value = 23624.428
$ws.Range("N104").Value = -30612.428
